$d = $word.ActiveDocument

# --- Edit 1: extend the first paragraph with a red "branch alternate" note ---
$p1 = $d.Paragraphs(1).Range
$p1.End = $p1.End - 1                 # exclude the paragraph mark
$p1.InsertAfter("  ")                 # two trailing spaces on the plain run

$s1 = $p1.End
$p1.InsertAfter("(This is a change – Ve")
$e1 = $p1.End
$d.Range($s1, $e1).Font.Color = 192   # 0xC00000 -> BGR long

$s2 = $p1.End
$p1.InsertAfter("rsion for branch alternate")
$e2 = $p1.End
$d.Range($s2, $e2).Font.Color = 192

$s3 = $p1.End
$p1.InsertAfter(")")
$e3 = $p1.End
$d.Range($s3, $e3).Font.Color = 192

# --- Edit 2: append a new, empty, shaded paragraph at the very end of the body ---
$endPos = $d.Content.End
$r = $d.Range($endPos, $endPos)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$r.InsertXML($xml)
